$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 10 (shifts old rows 10-12 down to 12-14)
$ws.Range("A10:R11").EntireRow.Insert()

# New row 10: Ciboulette Primera, date 2023-08-08 (45146)
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 45146
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112039
$ws.Cells.Item(10, 7).Value = "Ciboulette"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 2500
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2500
$ws.Cells.Item(10, 14).Value = "`$/docena de atados"
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 833
$ws.Cells.Item(10, 17).Value = 3
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# New row 11: Ciboulette Segunda, date 2023-08-08 (45146)
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 45146
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100112039
$ws.Cells.Item(11, 7).Value = "Ciboulette"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Segunda"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 2000
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 2000
$ws.Cells.Item(11, 14).Value = "`$/docena de atados"
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 667
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"
